$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A-C unchanged text, D/E get new headers) ---
$ws.Range("A1").Value = "Módulo"
$ws.Range("B1").Value = "Tema"
$ws.Range("C1").Value = "Cantidad de clases"

# --- Data rows (A=module #, B=topic, C=class count, D/E = "." placeholder) ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Inferencia causal"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "[Slide 1](https://drive.google.com/file/d/11S5sW3UUnEfDna5LPVoeTchjwbAmg1iw/view?usp=sharing)"

$ws.Range("D1").Value = "Presentaciones"
$ws.Range("E1").Value = "Laboratorios"

$ws.Range("E2").Value = "."

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Análisis bajo asignación aleatoria simple"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = "."
$ws.Range("E3").Value = "."

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Estrategias de aleatorización y análisis basado en diseño"
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = "."
$ws.Range("E4").Value = "."

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Diseños experimentales"
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = "."
$ws.Range("E5").Value = "."

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Análisis estadístico de experimentos"
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = "."
$ws.Range("E6").Value = "."

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Tipos de experimentos en ciencias sociales"
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = "."
$ws.Range("E7").Value = "."

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Pre-registro de experimentos"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = "."
$ws.Range("E8").Value = "."

# --- Column F is no longer used ---
$ws.Range("F1:F8").ClearContents()

# --- Selection / view matches the saved state in the target file ---
$ws.Range("G5").Select()
